# Add a new "N=200000" results sheet at the end of the workbook (after the
# existing "N=150000" sheet) and fill it with the bucket-sort timing data,
# mirroring the layout of the other "N=..." sheets.

$wb = $excel.ActiveWorkbook

# Insert the new worksheet right after the current last sheet so it ends up
# as the final tab, matching the author's sheet order.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "N=200000"

# Header row
$ws.Range("A1").Value = "Execução"
$ws.Range("B1").Value = "Tempo (ms)"

# Individual run measurements
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "242.8880 ms"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "237.9081 ms"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "236.9211 ms"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "241.9710 ms"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "271.4119 ms"

# Summary rows
$ws.Range("A7").Value = "Média"
$ws.Range("B7").Value = "246.2200 ms"

$ws.Range("A8").Value = "Desvio Padrão"
$ws.Range("B8").Value = "14.3121 ms"
